$wb = $excel.ActiveWorkbook

$wsInsumos = $wb.Worksheets.Item("Insumos")
$wsProcesados = $wb.Worksheets.Item("Procesados")
$wsNoProcesados = $wb.Worksheets.Item("NoProcesados")

# --- Insumos: clear row 5 (Allisson Flores Espinoza, duplicate) and part of row 6 (Eunice .. Hernandez) ---
$wsInsumos.Range("A5:C5").ClearContents()
$wsInsumos.Range("A6").ClearContents()
$wsInsumos.Range("C6").ClearContents()

# --- Procesados: first processed email ---
$wsProcesados.Range("A2").Value = "ivan.hernandez@beeckerco.com"

# --- NoProcesados: row for the duplicate-name record ---
$wsNoProcesados.Range("A2").Value = "Allisson"
$wsNoProcesados.Range("B2").Value = "Flores"
$wsNoProcesados.Range("C2").Value = "Espinoza"
$wsNoProcesados.Range("D2").Value = "Nombre repetido"

# --- NoProcesados: row for the missing-apellido-paterno record ---
$wsNoProcesados.Range("A3").ClearFormats()
$wsNoProcesados.Range("A3").Value = "Eunice"
$wsNoProcesados.Range("B3").Clear()
$wsNoProcesados.Range("C3").ClearFormats()
$wsNoProcesados.Range("C3").Value = "Hernández"
$wsNoProcesados.Range("D3").ClearFormats()
$wsNoProcesados.Range("D3").Value = "Sin apellido Paterno"

# --- Procesados: remaining processed emails ---
$wsProcesados.Range("A3").Value = "ivan.aparicio@beeckerco.com"
$wsProcesados.Range("A4").Value = "allisson.espinoza@beeckerco.com"
$wsProcesados.Rows.Item(4).RowHeight = 15.75

# --- View state: Procesados becomes the active/selected sheet ---
$wsNoProcesados.Range("G13").Select()
$wsProcesados.Activate()
$wsProcesados.Range("A6").Select()
